$d = $word.ActiveDocument

# --- Paragraph 1: keep the <w:br/> between the two runs, replace each half via scoped Find/Replace ---
$p1 = $d.Paragraphs.Item(1).Range
$p1.Find.Execute('המאמר היומי של מייק - 06.02.25', $true, $true, $false, $false, $false, $true, 1, $false, 'המאמר היומי של מייק - 05.02.25', 2) | Out-Null
$p1b = $d.Paragraphs.Item(1).Range
$p1b.Find.Execute('SMALL LANGUAGE MODELS: SURVEY, MEASUREMENTS, AND INSIGHTS', $true, $true, $false, $false, $false, $true, 1, $false, 'Deep Generative Models through the Lens of the Manifold Hypothesis: A Survey and New Connections', 2) | Out-Null

# --- Paragraphs 2..29: replace each paragraph's content with the new text. Paragraph 6 previously
#     carried a stray leading-space xml:space="preserve", so it goes through scoped Find/Replace
#     (which recomputes the attribute off the new text) instead of a bare Range.Text= assignment. ---
$d.Paragraphs.Item(2).Range.Text = 'תמצית המאמר: '
$d.Paragraphs.Item(3).Range.Text = 'רציתם לדעת למה מודלי דיפוזיה ניצחו את הגאנים, VAE וכל השאר מזווית מתמטית? רוצים להבין בעזרת מתמטיקה למה מודלי דיפוזיה לטנטיים עובדים מעולה? תצללו לסקירה הזו…'
$d.Paragraphs.Item(4).Range.Text = 'מאמר זה מציע חקירה מקיפה של מודלים גנרטיביים עמוקים (DGMs) תחת המסגרת של השערת היריעה, הטוענת שדאטה בעל ממד גבוה נמצאים לעתים קרובות על תת-יריעה בעלת ממד נמוך יותר המוטמעת בתוך המרחב המקורי (במאמר נקרא אמביינטי). המחברים מספקים הסבר מדוע מודלים כמו מודלי דיפוזיה ו- GANs מסוימות מציגים ביצועים טובים יותר מאחרים, כולל שיטות מבוססות נראות כמו אוטואנקודרים וריאציוניים (VAEs) וזרימות נורמליזציה (NFs). על ידי אימוץ נקודת מבט מבוססת יריעה, המחברים מספקים תובנות לגבי המגבלות המובנות של גישות קיימות תוך יצירת קשרים תיאורטיים חדשים בין DGMs והסעה אופטימלית..'
$d.Paragraphs.Item(5).Range.Text = 'המחקר בולט בכך שהוא מוכיח באופן פורמלי את חוסר היציבות הנומרית המובנית שמודלים מבוססי נראות בממד גבוה  חווים כאשר הם מנסים לייצג דאטה על יריעה, ומציע פרשנות חדשה של DGMs דו-שלביים כמקרבים של מרחק וסרשטיין בין התפלגות המודל להתפלגות הדאטה האמיתי.'
$pp6 = $d.Paragraphs.Item(6).Range
$pp6.Find.Execute(' 1. חידושים ארכיטקטוניים מאמר המחקר בוחן לעומק את האלמנטים הארכיטקטוניים המבדילים בין SLM ל-LLM, תוך הדגשת השינויים המשפרים את היעילות במכשירים עם משאבים מוגבלים.', $true, $true, $false, $false, $false, $true, 1, $false, 'נקודות מרכזיות', 2) | Out-Null
$d.Paragraphs.Item(7).Range.Text = '1. סקירה של מודלים DGM מודעי-יריעה ולא-מודעי-יריעה (manifold-aware and manifold unaware)'
$d.Paragraphs.Item(8).Range.Text = 'מודלים לא-מודעי-יריעה: מודלים אלה אינם מתחשבים באופן מפורש במבנה היריעה של דאטה. דוגמאות כוללות VAEs, NFs ומודלים מבוססי אנרגיה. מודלים כאלה נוטים להתאמת יתר ליריעה, כאשר הצפיפויות שואפות לאינסוף לאורך היריעה אך נכשלים בשעורכה של ההתפלגות בתוכה.'
$d.Paragraphs.Item(9).Range.Text = 'מודלים מודעי-יריעה: מודלים אלה מוסיפים רעש כדי לפזר את מסת ההסתברות מעבר ליריעה או מאפטמים פונקציות יעד שאינם מגבילות את ההתפלגות על היריעה שתופסות באופן לא מפורש את מבנה היריעה. דוגמאות כוללות מודלי דיפוזיה, התאמת זרימה מותנית (conditional flow models), ו-Wasserstein GANs.'
$d.Paragraphs.Item(10).Range.Text = '2. חוסר יציבות נומרית של שיטות מבוססות נראות'
$d.Paragraphs.Item(11).Range.Text = 'אחת התרומות התיאורטיות המרכזיות היא ההוכחה שמודלים מבוססי נראות סובלים מחוסר יציבות מספרית בלתי נמנע כאשר הם מנסים למדל הדאטה הנתמך על יריעה. המחברים מדגימים שכאשר צפיפויות המודל מנסות להתרכז על היריעה, פונקציית הנראות הופכת לבלתי מוגבלת, מה שמוביל לפתרונות מנוונים(זה קורה הרבה ב-VAE וב-GANים רגילים).'
$d.Paragraphs.Item(12).Range.Text = 'מתמטית, אם P_X  התפלגות הדאטה ב- Rd בעלת תומך של יריעה M בעלת ממד פנימי d*< d, עבור כל סדרה של מודלים מבוססי נראות { P_{X, θ_t} המקרבים את התפלגות דאטה מתקיים: '
$d.Paragraphs.Item(13).Range.Text = 'תוצאה זו מרמזת שצפיפויות במרחב הדאטה מתבדרות באופן מובנה כאשר הן מנסות למדל התפלגויות הנתמכות על יריעה, מה שהופך את היעדים מבוססי הנראות לבעייתיים עבור דאטה כזה כאלה.'
$d.Paragraphs.Item(14).Range.Text = '3. מגבלות מרחק KL:'
$d.Paragraphs.Item(15).Range.Text = 'המחברים מדגישים שמרחק KL, יעד נפוץ לאימון DGMs, הופך ללא יעיל בלמידת היריעה. הבעיה העיקרית מתעוררת כי KL מניח ששתי ההתפלגויות חולקות את אותה תומך (support). אולם כאשר משווים צפיפות של מודל במרחב הדאטה {p_{X,θ   עם התפלגות דאטה P_ הנתמכת על יריעה, ה-KL הופך לאינסופי:'
$d.Paragraphs.Item(16).Range.Text = 'תופעה זו מתרחשת כי P_x מקצה הסתברות שאינה אפס רק לנקודות על היריעה, בעוד {p_{X,θ מפזר מסת הסתברות על פני כל המרחב האמביינטי. כתוצאה מכך, מקסום הנראות, השקול למזעור את ה-KL, נכשל במתן אות למידה משמעותי.'
$d.Paragraphs.Item(17).Range.Text = '4. מרחק וסרשטיין כיעד חלופי'
$d.Paragraphs.Item(18).Range.Text = 'כדי להתמודד עם מגבלות ה-KL, המחברים מקדמים את השימוש במרחקי וסרשטיין(זה עובד לא רע בגאנים כאמור), שנשארים מוגדרים היטב גם כאשר להתפלגויות יש תמיכות לא תואמות. מרחק וסרשטיין-1 בין התפלגויות p ו-q מוגדר כ:'
$d.Paragraphs.Item(19).Range.Text = 'כאשר (Π(p,q מסמן את קבוצת ההתפלגויות המשותפות עם בעלות התפלגות שולית p ו-q. בניגוד ל- KL, מרחק וסרשטיין ממטר התכנסות חלשה, מה שהופך אותו ליעד חסין לאימון DGMs בתרחיש היריעה.'
$d.Paragraphs.Item(20).Range.Text = '5. פרשנות של מודלים לטנטיים'
$d.Paragraphs.Item(21).Range.Text = 'המחברים מספקים פרשנות חדשה של DGMs לטנטיים שקודם לומדים ייצוג בממד נמוך של יריעת דאטה ואז ממדלים את ההתפלגות בתוך ייצוג זה. הם מראים שמודלים אלה ממזערים באופן יעיל חסם עליון של מרחק וסרשטיין בין התפלגות המודל להתפלגות דאטה האמיתית:'
$d.Paragraphs.Item(22).Range.Text = 'כאשר שגיאת השחזור מודדת עד כמה טוב היריעה שנלמדה מקרבת את יריעת הדאטה האמיתית, והמרחק בין ההתפלגות מכמת את ההבדל בין התפלגויות בתוך היריעה שנלמדה. תוצאה זו מספקת הצדקה תיאורטית להצלחה האמפירית של מודלי דיפוזיה לטנטיים וגישות דו-שלביות אחרות.'
$d.Paragraphs.Item(23).Range.Text = 'תובנות מתמטיות'
$d.Paragraphs.Item(24).Range.Text = 'משפט חוסר היציבות המספרית'
$d.Paragraphs.Item(25).Range.Text = 'המחברים מוכיחים באופן פורמלי שעבור כל התפלגות הדאטה P הנתמכת על יריעה וכל סדרה של צפיפויות מודל במימד הדאטה Qₙ, פונקציית יעד בצורה של נראות לא מתכנסת מקסימום. תוצאה זו נגזרת מניתוח התנהגות הצפיפויות על יריעות בממד נמוך תוך שימוש בתכונות של גיאומטריה דיפרנציאלית ותורת המידה (די כבד האמת).'
$d.Paragraphs.Item(26).Range.Text = 'מזעור מרחק וסרשטיין:'
$d.Paragraphs.Item(27).Range.Text = 'על ידי הצגת מודלים דו-שלביים כמקרבים של מרחק וסרשטיין, המחברים מבססים קשר בין למידת יריעה וטרנספורט אופטימלי. תובנה זו לא רק מסבירה את הביצועים העדיפים של מודלי דיפוזיה לטנטיים אלא גם מספקת מסגרת עקרונית לתכנון DGMs חדשים.'
$d.Paragraphs.Item(28).Range.Text = 'הסבר קריסת מודים:'
$d.Paragraphs.Item(29).Range.Text = 'המחברים מראים שקריסת מודים ב-VAEs ו-GANs ניתנת להבנה כתוצאה של התאמת יתר ליריעה, כאשר צפיפויות המודל מתבדרות לאורך תתי-קבוצות של היריעה מבלי ״לתפוס״ את התפלגות הדאטה האמיתית.'

# --- Append five brand-new paragraphs after paragraph 29 ---
$tailTexts = @('מודלי דיפוזיה:', 'ההצלחה של מודלי דיפוזיה מיוחסת ליכולתם להתחשב באופן מרומז במבנה היריעה על ידי פיזור מסת הסתברות מעבר ליריעה. המחברים מספקים ניתוח מפורט של מודלי דיפוזיה מבוססי-ציון וגרסאות חבויות שלהם.', 'סיכום', 'מאמר זה מספק חקירה קפדנית ומעמיקה של DGMs דרך עדשת השערת היריעה. על ידי זיהוי המגבלות של שיטות מבוססות-נראות והדגשת היתרונות של מרחקי וסרשטיין ומודלים לטנטיים, המחברים סוללים את הדרך לפיתוח מודלים גנרטיביים יעילים יותר. ', 'https://arxiv.org/abs/2404.02954')
foreach ($t in $tailTexts) {
    $d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()
    $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = $t
}

Write-Output ("FinalParaCount=" + $d.Paragraphs.Count)